$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "320.30") need to be
# forced to Text format first, otherwise Excel auto-converts them to a number
# (dropping the trailing zero) instead of keeping the literal display string.
$textCells = @(
    "D5",
    "D6",
    "D9",
    "D10",
    "D11",
    "D14",
    "D15",
    "D17",
    "D19",
    "D21",
    "D22",
    "D23",
    "D27",
    "D29",
    "D30",
    "D31",
    "D34",
    "D36",
    "D45",
    "D46",
    "D47",
    "D49",
    "D50",
    "D51",
)
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D5").Value = '320.30'
$ws.Range("D6").Value = '92.14'
$ws.Range("D9").Value = '0.510'
$ws.Range("D10").Value = '32.96'
$ws.Range("D11").Value = '0.0851'
$ws.Range("D14").Value = '6.88'
$ws.Range("D15").Value = '15.45'
$ws.Range("D17").Value = '0.791'
$ws.Range("D19").Value = '6.44'
$ws.Range("D21").Value = '71.40'
$ws.Range("D22").Value = '11.23'
$ws.Range("D23").Value = '239.50'
$ws.Range("D27").Value = '24.90'
$ws.Range("D29").Value = '9.73'
$ws.Range("D30").Value = '36.49'
$ws.Range("D31").Value = '157.07'
$ws.Range("D34").Value = '0.0765'
$ws.Range("D36").Value = '17.12'
$ws.Range("D45").Value = '18.65'
$ws.Range("D46").Value = '2.96'
$ws.Range("D47").Value = '9.45'
$ws.Range("D49").Value = '97.47'
$ws.Range("D50").Value = '75.67'
$ws.Range("D51").Value = '66.88'

# Restore the default (unstyled) cell style now that the value is locked in as text.
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).Style = "Normal"
}

# Remaining cells (percentages and non-ambiguous text) can be set directly.
$ws.Range("D2").Value = '41.703.51'
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("D3").Value = '2.472.60'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  +1.26%  '
$ws.Range("E6").Value = '  -0.92%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -1.02%  '
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("D13").Value = '2.853.22'
$ws.Range("E13").Value = '  +0.13%  '
$ws.Range("E14").Value = '  -0.30%  '
$ws.Range("D16").Value = '2.461.22'
$ws.Range("E17").Value = '  +1.36%  '
$ws.Range("D18").Value = '41.630.59'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("D20").Value = '0.0₃0938'
$ws.Range("E20").Value = '  -1.26%  '
$ws.Range("E21").Value = '  +0.41%  '
$ws.Range("E22").Value = '  -2.80%  '
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("E24").Value = '  +1.09%  '
$ws.Range("E25").Value = '  +0.90%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("E28").Value = '  -1.07%  '
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("E30").Value = '  +1.60%  '
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("E32").Value = '  -1.82%  '
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("E35").Value = '  -0.96%  '
$ws.Range("E36").Value = '  -2.68%  '
$ws.Range("E37").Value = '  +0.94%  '
$ws.Range("E38").Value = '  -0.45%  '
$ws.Range("E39").Value = '  +1.35%  '
$ws.Range("E40").Value = '  -0.33%  '
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("E42").Value = '  -2.71%  '
$ws.Range("D43").Value = '2.003.88'
$ws.Range("E43").Value = '  +1.57%  '
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("E46").Value = '  +0.21%  '
$ws.Range("E47").Value = '  +4.44%  '
$ws.Range("D48").Value = '2.730.62'
$ws.Range("E48").Value = '  +1.02%  '
$ws.Range("E49").Value = '  +0.61%  '
$ws.Range("E50").Value = '  +3.71%  '
$ws.Range("E51").Value = '  -0.61%  '
